# SO_Test.xlsx - add new completed order as row 24 on Sheet1, following the
# same column layout used by the existing rows (A:AD).
#
# Columns with date-looking text ("03/11/2024") must be stored as literal
# text, not auto-converted to Excel date serials, so each such cell is
# switched to the Text number format before the value is written (and the
# format is cleared again afterwards so no extra style lingers on the cell).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$row = 24

function Set-TextValue($r, $c, $val) {
    $cell = $ws.Cells.Item($r, $c)
    $cell.NumberFormat = "@"
    $cell.Value = $val
    $cell.ClearFormats()
}

function Set-BlankCell($r, $c) {
    $cell = $ws.Cells.Item($r, $c)
    $cell.NumberFormat = "@"
    $cell.Value = ""
    $cell.ClearFormats()
}

Set-TextValue $row 1  "03/11/2024"
Set-TextValue $row 2  "SO240311007"
Set-TextValue $row 3  "YES"
Set-TextValue $row 4  "ab"
Set-TextValue $row 5  "asdfjkhl"
Set-TextValue $row 6  "(789)456-1230"
Set-TextValue $row 7  "YES"
Set-TextValue $row 8  "YES"
Set-TextValue $row 9  "aszdg"
Set-TextValue $row 10 "adfsg"

$ws.Cells.Item($row, 11).Value = 50
$ws.Cells.Item($row, 12).Value = 900

Set-TextValue $row 13 " "
Set-TextValue $row 14 "ytouse"
Set-TextValue $row 15 "DVD"
Set-TextValue $row 16 "abake"
Set-TextValue $row 17 "PICKUP"

Set-BlankCell $row 18
Set-BlankCell $row 19
Set-BlankCell $row 20
Set-BlankCell $row 21

Set-TextValue $row 22 "YES"
Set-TextValue $row 23 "03/11/2024"
Set-TextValue $row 24 "abake"
Set-TextValue $row 25 "YES"
Set-TextValue $row 26 "03/11/2024"
Set-TextValue $row 27 "abake"
Set-TextValue $row 28 "YES"
Set-TextValue $row 29 "03/11/2024"
Set-TextValue $row 30 "abake"
